$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29..126 down to 30..127.
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with its data (a new Mango price record).
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44672
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100108
$ws.Range("H29").Value = "Tropicales y subtropicales"
$ws.Range("I29").Value = 100108002
$ws.Range("J29").Value = "Mango"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 7000
$ws.Range("O29").Value = 7000
$ws.Range("P29").Value = 7000
$ws.Range("Q29").Value = "$/bandeja 4 kilos"
$ws.Range("R29").Value = "Perú"
$ws.Range("S29").Value = 1750
$ws.Range("T29").Value = 4
